# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Rebuilds the detail rows (16-35) of the "Estado de Cuenta" sheet so that
# MARISELA ARRIETA DOMINGUEZ (CC 34948070) and LUZ MARINA ARRIETA DOMINGUEZ
# (CC 34948055) alternate, one row each, across the full set of periods
# (2004 .. 2012, 2101), and refreshes the mora / salary values for the
# updated periods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$worker1Doc = "34948070"
$worker1Name = "MARISELA ARRIETA DOMINGUEZ"
$worker2Doc = "34948055"
$worker2Name = "LUZ MARINA ARRIETA DOMINGUEZ"

$periods = @("2004","2005","2006","2007","2008","2009","2010","2011","2012","2101")

$row = 16
foreach ($periodo in $periods) {
    if ($periodo -eq "2101") {
        $valorMora = 26500
    } else {
        $valorMora = 33125
    }
    $salario = 828116

    $ws.Range("B" + $row).Value = "CC"
    $ws.Range("C" + $row).Value = $worker1Doc
    $ws.Range("D" + $row).Value = $worker1Name
    $ws.Range("E" + $row).Value = $periodo
    $ws.Range("F" + $row).Value = $valorMora
    $ws.Range("G" + $row).Value = $salario
    $row = $row + 1

    $ws.Range("B" + $row).Value = "CC"
    $ws.Range("C" + $row).Value = $worker2Doc
    $ws.Range("D" + $row).Value = $worker2Name
    $ws.Range("E" + $row).Value = $periodo
    $ws.Range("F" + $row).Value = $valorMora
    $ws.Range("G" + $row).Value = $salario
    $row = $row + 1
}
